# Adds two new worksheets - "CaseDetailStat" and "CaseDetailStat_Message" -
# to the workbook, populating them with the "case detail" file listing and
# the accompanying run-log/message content, mirroring the existing
# StatOutput / StatOutput_Message sheet pair.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Create the two new sheets, appended after the last existing sheet, in
# the same relative order as StatOutput -> StatOutput_Message.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$caseDetailStat = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$caseDetailStat.Name = "CaseDetailStat"

$caseDetailStatMsg = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $caseDetailStat)
$caseDetailStatMsg.Name = "CaseDetailStat_Message"

# ---------------------------------------------------------------------
# CaseDetailStat: header row + single data row describing one file
# associated with case NCATS-COP01CCB010072.
# ---------------------------------------------------------------------
$caseDetailStat.Range("A1").Value = "File Name"
$caseDetailStat.Range("B1").Value = "File Type"
$caseDetailStat.Range("C1").Value = "Association"
$caseDetailStat.Range("D1").Value = "Description"
$caseDetailStat.Range("E1").Value = "Format"
$caseDetailStat.Range("F1").Value = "Size"

$caseDetailStat.Range("A2").Value = "CCB010072.pdf"
$caseDetailStat.Range("B2").Value = "Pathology Report"
$caseDetailStat.Range("C2").Value = "diagnosis"
# Description is blank for this file (Neo4j returned an empty string).
$caseDetailStat.Range("D2").Value = ""
$caseDetailStat.Range("E2").Value = "pdf"

# Keep the file-size value textual ("57.732421875") instead of letting
# Excel auto-coerce it to a number.
$caseDetailStat.Range("F2").NumberFormat = "@"
$caseDetailStat.Range("F2").Value = "57.732421875"

# ---------------------------------------------------------------------
# CaseDetailStat_Message: connection/run log, repeated for the three
# queries that were executed (case detail query, breed-count stat query,
# and the new case-detail-stat count query), matching the pattern already
# used on the Message / CypherOutput_Message / StatOutput_Message sheets.
# ---------------------------------------------------------------------
$neo4jUrl = "bolt://ncias-q2251-c.nci.nih.gov:7687"
$outputPath = "C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC02_Canine_Filter_Breed-AmerStaffd_Neo4jData.xlsx"

$caseDetailsQuery = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN ['American Staffordshire Terrier'] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'') AS ``Case ID`` , coalesce(s.clinical_study_designation,'') AS ``Study Code`` , coalesce(s.clinical_study_type,'') AS  ``Study Type``, coalesce(demo.breed,'') AS Breed , coalesce(diag.disease_term,'') AS Diagnosis , coalesce(diag.stage_of_disease,'') AS ``Stage of Disease`` ,  coalesce(demo.patient_age_at_enrollment,'') AS Age , coalesce(demo.sex,'') AS Sex , coalesce(demo.neutered_indicator,'') AS  ``Neutered Status``"

$breedCountQuery = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['American Staffordshire Terrier']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

$caseDetailStatQuery = $breedCountQuery

$lines = @(
    "Neo4j_URL:", $neo4jUrl, "User_name:", "neo4j", "PWD:", "icdcDBneo4j0", "Cypher:", $caseDetailsQuery, "Output:", $outputPath,
    "Neo4j_URL:", $neo4jUrl, "User_name:", "neo4j", "PWD:", "icdcDBneo4j0", "Cypher:", $breedCountQuery, "Output:", $outputPath,
    "Neo4j_URL:", $neo4jUrl, "User_name:", "neo4j", "PWD:", "icdcDBneo4j0", "Cypher:", $caseDetailStatQuery, "Output:", $outputPath
)

for ($i = 0; $i -lt $lines.Length; $i++) {
    $row = $i + 1
    $caseDetailStatMsg.Cells.Item($row, 1).Value = $lines[$i]
}

# Restore the originally-selected/active sheet.
$wb.Worksheets.Item(1).Activate()
